$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates whose day-of-month is <=12 are ambiguous (Excel could read them as
# MM-DD-YYYY) and get silently auto-converted into a date serial number on
# plain assignment. Force those through a text number format, then restore
# the cell to the (unstyled) Normal style once the literal text is locked
# in, so no explicit style index lingers on the cell.
function Set-TextDate($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("A3").Value = "28-07-2022"
Set-TextDate "A4" "01-08-2022"
Set-TextDate "A5" "04-08-2022"
Set-TextDate "A6" "08-08-2022"
Set-TextDate "A7" "11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
Set-TextDate "A13" "01-09-2022"
Set-TextDate "A14" "05-09-2022"
Set-TextDate "A15" "08-09-2022"
Set-TextDate "A16" "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D7").Value = 1
$ws.Range("G7").Value = 1

$ws.Range("D11").Value = 1
$ws.Range("G11").Value = 1

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0
